{"js": "// The document had a sentence split into 3 separate runs, with Word\n// grammar-checker <w:proofErr> markers (gramStart/gramEnd) wrapped around\n// the word \"issued\" in the middle:\n//\n//   \"Only controlled copies of this Manual are \" + [gramStart] \"issued\" [gramEnd] +\n//   \" and each will bear a unique number and be assigned to an individual.\"\n//\n// The edit consolidates this into a single run / single text node with no\n// proofErr markers:\n//\n//   \"Only controlled copies of this Manual are issued and each will bear a\n//    unique number and be assigned to an individual.\"\n//\n// Find the paragraph containing that sentence and rewrite its text in one\n// shot; replacing the whole paragraph (rather than a sub-range) collapses\n// the runs that used to be split apart by the grammar-check proofing marks\n// into a single clean run, matching how Word re-serializes the paragraph\n// once the sentence is no longer broken up.\nconst targetText =\n  \"Only controlled copies of this Manual are issued and each will bear a unique number and be assigned to an individual.\";\nconst searchText = \"Only controlled copies of this Manual are\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.indexOf(searchText) !== -1) {\n    para.insertText(targetText, Word.InsertLocation.replace);\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The paragraph \"Only controlled copies of this Manual are issued and each\n# will bear a unique number and be assigned to an individual.\" used to be\n# split across three separate runs, with Word's grammar-checker proofing\n# marks (<w:proofErr w:type=\"gramStart\"/> ... <w:proofErr w:type=\"gramEnd\"/>)\n# wrapped around the word \"issued\" in the middle:\n#\n#   \"Only controlled copies of this Manual are \" + [gramStart]\"issued\"[gramEnd] +\n#   \" and each will bear a unique number and be assigned to an individual.\"\n#\n# Running this sentence through Find & Replace (replacing it with its own,\n# identical text) makes Word re-serialize the paragraph as a single clean\n# run with no leftover proofing marks - exactly the edit recorded in the\n# diff.\n\n$d = $word.ActiveDocument\n\n$sentence = \"Only controlled copies of this Manual are issued and each will bear a unique number and be assigned to an individual.\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = $sentence\n$find.Replacement.Text = $sentence\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute($sentence, $false, $true, $false, $false, $false, $true, 1, $false, $sentence, 2) | Out-Null\n"}
